$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 294, pushing the existing data (rows 294-331) down to rows 297-334
$ws.Range("A294:A296").EntireRow.Insert()

# Copy formatting (style) of column D from the row right below (row 297, which used to be row 294)
# so that the new date cells keep the date style (s="2").
$ws.Range("D297").Copy()
$ws.Range("D294:D296").PasteSpecial(-4122) # xlPasteFormats

# Row 294 (new): Primera, Vol 300, Min 12000, Max 13000, Avg 12400, $/bandeja 7 kilos, Provincia de Melipilla, 1771
$ws.Cells.Item(294, 1).Value = 10
$ws.Cells.Item(294, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(294, 3).Value = "La Araucanía"
$ws.Cells.Item(294, 4).Value = 45212
$ws.Cells.Item(294, 5).Value = 9
$ws.Cells.Item(294, 6).Value = "Fruta"
$ws.Cells.Item(294, 7).Value = 100101
$ws.Cells.Item(294, 8).Value = "Berries"
$ws.Cells.Item(294, 9).Value = 100112025
$ws.Cells.Item(294, 10).Value = "Frutilla"
$ws.Cells.Item(294, 11).Value = "Sin especificar"
$ws.Cells.Item(294, 12).Value = "Primera"
$ws.Cells.Item(294, 13).Value = 300
$ws.Cells.Item(294, 14).Value = 12000
$ws.Cells.Item(294, 15).Value = 13000
$ws.Cells.Item(294, 16).Value = 12400
$ws.Cells.Item(294, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(294, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(294, 19).Value = 1771
$ws.Cells.Item(294, 20).Value = 7

# Row 295 (new): Segunda, Vol 100, Min 10000, Max 10000, Avg 10000, $/bandeja 7 kilos, Provincia de Melipilla, 1429
$ws.Cells.Item(295, 1).Value = 10
$ws.Cells.Item(295, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(295, 3).Value = "La Araucanía"
$ws.Cells.Item(295, 4).Value = 45212
$ws.Cells.Item(295, 5).Value = 9
$ws.Cells.Item(295, 6).Value = "Fruta"
$ws.Cells.Item(295, 7).Value = 100101
$ws.Cells.Item(295, 8).Value = "Berries"
$ws.Cells.Item(295, 9).Value = 100112025
$ws.Cells.Item(295, 10).Value = "Frutilla"
$ws.Cells.Item(295, 11).Value = "Sin especificar"
$ws.Cells.Item(295, 12).Value = "Segunda"
$ws.Cells.Item(295, 13).Value = 100
$ws.Cells.Item(295, 14).Value = 10000
$ws.Cells.Item(295, 15).Value = 10000
$ws.Cells.Item(295, 16).Value = 10000
$ws.Cells.Item(295, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(295, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(295, 19).Value = 1429
$ws.Cells.Item(295, 20).Value = 7

# Row 296 (new): Tercera, Vol 140, Min 7000, Max 7000, Avg 7000, $/bandeja 7 kilos, Provincia de Melipilla, 1000
$ws.Cells.Item(296, 1).Value = 10
$ws.Cells.Item(296, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(296, 3).Value = "La Araucanía"
$ws.Cells.Item(296, 4).Value = 45212
$ws.Cells.Item(296, 5).Value = 9
$ws.Cells.Item(296, 6).Value = "Fruta"
$ws.Cells.Item(296, 7).Value = 100101
$ws.Cells.Item(296, 8).Value = "Berries"
$ws.Cells.Item(296, 9).Value = 100112025
$ws.Cells.Item(296, 10).Value = "Frutilla"
$ws.Cells.Item(296, 11).Value = "Sin especificar"
$ws.Cells.Item(296, 12).Value = "Tercera"
$ws.Cells.Item(296, 13).Value = 140
$ws.Cells.Item(296, 14).Value = 7000
$ws.Cells.Item(296, 15).Value = 7000
$ws.Cells.Item(296, 16).Value = 7000
$ws.Cells.Item(296, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(296, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(296, 19).Value = 1000
$ws.Cells.Item(296, 20).Value = 7
